# Implementacion final de notificacion por correo (Enviado, Observado)
#
# Adds a new column F ("Enviado"/"Observado" notification flag) to the
# sheet, populating every data row (1-22) with the literal text "null"
# (a placeholder meaning "not yet notified"). This mirrors the author's
# change: a new shared string "null" is introduced and written into
# F1:F22, extending the used range from A1:E22 to A1:F22 and moving the
# active selection onto the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the new column F for every existing row with the "null" placeholder.
$ws.Range("F1:F22").Value = "null"

# Match the author's final selection: the whole new column is selected,
# anchored on F1.
$ws.Range("F1:F22").Select()
